$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39 - this shifts existing rows 39:100 down to 40:101
# and grows the used range from A1:T100 to A1:T101, matching the target diff.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new record's data.
$ws.Range("A39").Value = 9
$ws.Range("B39").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = 44662
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100101
$ws.Range("H39").Value = "Berries"
$ws.Range("I39").Value = 100101004
$ws.Range("J39").Value = "Frambuesa"
$ws.Range("K39").Value = "Sin especificar"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 220
$ws.Range("N39").Value = 8000
$ws.Range("O39").Value = 8000
$ws.Range("P39").Value = 8000
$ws.Range("Q39").Value = "`$/bandeja 2 kilos"
$ws.Range("R39").Value = "Provincia de Linares"
$ws.Range("S39").Value = 4000
$ws.Range("T39").Value = 2
